$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the numeric-looking price cells to remain plain text (as scraped),
# matching the workbooks existing inline-string data model before assigning values.
$textCells = @("D5", "D8", "D11", "D17", "D18", "D20", "D22", "D25", "D34", "D35", "D39", "D40", "D41", "D43", "D44", "D45", "D48", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated crypto price/volume data as scraped on Mon Oct 16 04:29:38 UTC 2023
$ws.Range("D2").Value = "27.255.56"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "1.563.71"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  -0.49%  "
$ws.Range("D5").Value = "211.02"
$ws.Range("E5").Value = "  +1.26%  "
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("D8").Value = "22.19"
$ws.Range("E8").Value = "  +0.63%  "
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("E10").Value = "  -0.52%  "
$ws.Range("D11").Value = "0.0871"
$ws.Range("E11").Value = "  +1.84%  "
$ws.Range("D12").Value = "1.787.47"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").Value = "1.563.17"
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("D16").Value = "27.266.37"
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("D17").Value = "61.84"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "218.05"
$ws.Range("E18").Value = "  +0.62%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0703"
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "7.44"
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("E21").Value = "  -0.44%  "
$ws.Range("D22").Value = "4.14"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("E23").Value = "  +1.89%  "
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").Value = "151.52"
$ws.Range("E25").Value = "  -1.00%  "
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("E27").Value = "  +1.02%  "
$ws.Range("E28").Value = "  -0.56%  "
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("E30").Value = "  +2.06%  "
$ws.Range("E31").Value = "  -0.68%  "
$ws.Range("B33").Value = "Maker"
$ws.Range("C33").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D33").Value = "1.457.72"
$ws.Range("E33").Value = "  +2.06%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "3.18"
$ws.Range("E34").Value = "  +0.45%  "
$ws.Range("D35").Value = "1.11"
$ws.Range("E35").Value = "  +5.26%  "
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").Value = "0.540"
$ws.Range("E39").Value = "  +1.06%  "
$ws.Range("D40").Value = "5.86"
$ws.Range("E40").Value = "  -0.44%  "
$ws.Range("D41").Value = "0.814"
$ws.Range("E41").Value = "  +0.55%  "
$ws.Range("E42").Value = "  -0.42%  "
$ws.Range("D43").Value = "2.34"
$ws.Range("E43").Value = "  +1.22%  "
$ws.Range("D44").Value = "0.978"
$ws.Range("E44").Value = "  -2.32%  "
$ws.Range("D45").Value = "64.41"
$ws.Range("E45").Value = "  -0.42%  "
$ws.Range("E46").Value = "  +0.50%  "
$ws.Range("D47").Value = "1.701.40"
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("D48").Value = "85.86"
$ws.Range("E48").Value = "  -1.22%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("E50").Value = "  +1.00%  "
$ws.Range("D51").Value = "0.0948"
$ws.Range("E51").Value = "  -1.18%  "
